$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) to reflect new field names
$ws.Range("C1").Value = "x_array"
$ws.Range("D1").Value = "y_array"
$ws.Range("E1").Value = "i_sense"
$ws.Range("F1").Value = "entx"
$ws.Range("G1").Value = "enty"

# Remove now-unused columns H and I entirely, shifting remaining cells left
$ws.Range("H1:I3").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# Delete row 3 entirely, shifting remaining rows up
$ws.Rows.Item(3).Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
